# Generate Report for Handoff
# Updates the "b.md" file row across the Overview, zh-cn and de-de sheets
# to reflect that it is now "Ready for handoff" with a new handoff file
# (b.63290e5768f688058c7b37413b0a5c26c308f864.*) and timestamp.

$wb = $excel.ActiveWorkbook

# --- Overview sheet (row for b.md is row 3) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-03-22 12:38:48"

# --- zh-cn sheet (row for b.md is row 3) ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZh.Range("E3").Value = "2016-03-22 12:38:44"

foreach ($h in $wsZh.Hyperlinks) {
    if ($h.Range.Address() -eq '$D$3') {
        $h.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
    }
}

# --- de-de sheet (row for b.md is row 3) ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDe.Range("E3").Value = "2016-03-22 12:38:48"

foreach ($h in $wsDe.Hyperlinks) {
    if ($h.Range.Address() -eq '$D$3') {
        $h.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
    }
}
